# Fixed internal temp conversion
# "Internal" sheet: the thermistor temperature formula in column B had the
# numerator of the division backwards (0.7012 - x) instead of (x - 0.7012),
# and the "approx" formula in column C used the wrong linear model
# (2*A-410) instead of the correct one ((230-A)*2). Re-enter both formulas
# (anchor cell + the shared-formula fill range) so every dependent cell
# (column D, and the chart built on B:C) recalculates with the corrected
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Internal")

# --- Column B: corrected thermistor conversion formula ---
$ws.Range("B2").Formula = "=IF((3.3*A2/1023)>0.7012,25-((3.3*A2/1023)-0.7012)/0.001769,25-((3.3*A2/1023)-0.7012)/0.001646)"
$ws.Range("B3:B37").Formula = "=IF((3.3*A3/1023)>0.7012,25-((3.3*A3/1023)-0.7012)/0.001769,25-((3.3*A3/1023)-0.7012)/0.001646)"

# --- Column C: corrected linear approximation formula ---
$ws.Range("C2").Formula = "=(230-A2)*2"
$ws.Range("C3:C37").Formula = "=(230-A3)*2"

$wb.Application.Calculate()
